$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that were removed entirely (naive forecaster bug fix: bad data points dropped)
$ws.Range("C2").Value = $null
$ws.Range("E2").Value = $null
$ws.Range("C3").Value = $null
$ws.Range("C4").Value = $null

# Update recalculated values (floating point precision shifts from corrected averaging)
$ws.Range("E4").Value = 13.08276537368063
$ws.Range("C5").Value = -8.992252553594259
$ws.Range("E5").Value = -19.76480035196673
$ws.Range("C6").Value = -7.266312015249799
$ws.Range("E6").Value = 12.31225042954256
$ws.Range("C7").Value = 7.007132997505217
$ws.Range("C9").Value = 8.866443976147087
$ws.Range("C10").Value = 9.469137444079955
$ws.Range("E10").Value = 10.69920649119718
$ws.Range("C11").Value = 3.0013062146236
$ws.Range("C12").Value = 3.358206407534969
$ws.Range("E12").Value = 4.390489499870132
$ws.Range("C13").Value = -2.90476933598719
$ws.Range("E14").Value = 3.502435351035582
$ws.Range("E15").Value = 9.131012060398703
$ws.Range("C16").Value = 3.901355411819685
$ws.Range("E16").Value = 6.143002545701304
$ws.Range("E18").Value = 4.555278923792594
$ws.Range("E20").Value = 0.9515943257393467
$ws.Range("C21").Value = 4.073887526082043
$ws.Range("C22").Value = 5.246209615995689
$ws.Range("C23").Value = 7.340964210079837
$ws.Range("E23").Value = 8.197760099691219
$ws.Range("E24").Value = -4.308894244053663
$ws.Range("C25").Value = 5.152630504861988
$ws.Range("E25").Value = 6.518301903862955
$ws.Range("C26").Value = 4.862559663742938
$ws.Range("C27").Value = 3.497157880977597
$ws.Range("E29").Value = -0.7915059299106075
$ws.Range("C30").Value = 2.76474001115945
$ws.Range("E30").Value = 0.3611963426345843
$ws.Range("C32").Value = -2.305533699949858
$ws.Range("C33").Value = -12.52375957300176
$ws.Range("E33").Value = -47.88137131999108
$ws.Range("C34").Value = -7.260793671746447
$ws.Range("E34").Value = 21.21858006100774
$ws.Range("C35").Value = 6.393774768527805
$ws.Range("E35").Value = 15.67089147385301
$ws.Range("E36").Value = 16.14645080511215
$ws.Range("C38").Value = 4.097586525396246
$ws.Range("C41").Value = 7.775484240847486
$ws.Range("C42").Value = 7.824284864703768
$ws.Range("C43").Value = -0.7754456294374767
$ws.Range("E43").Value = -5.245870741613645
$ws.Range("E44").Value = -3.570724939213787
$ws.Range("C45").Value = -0.4725309217548324
$ws.Range("C46").Value = -1.245022353133318
$ws.Range("C47").Value = -5.658719874383722
$ws.Range("E47").Value = -6.559334139663897
$ws.Range("C49").Value = -2.497020900647684
$ws.Range("E49").Value = 0.3065176867123487
$ws.Range("C53").Value = 3.002653050584203
